$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("L4").Value = 3
$ws.Range("M4").Value = 3
$ws.Range("N4").Value = 6
$ws.Range("O4").Value = 7
$ws.Range("P4").Value = 8
$ws.Range("S4").Value = 8
$ws.Range("T4").Value = 8
$ws.Range("U4").Value = 9
$ws.Range("V4").Value = 7
$ws.Range("W4").Value = 6
$ws.Range("X4").Value = 2

# Row 5
$ws.Range("M5").Value = 3
$ws.Range("N5").Value = 5
$ws.Range("O5").Value = 6
$ws.Range("P5").Value = 6
$ws.Range("Q5").Value = 6
$ws.Range("R5").Value = 7
$ws.Range("S5").Value = 7
$ws.Range("T5").Value = 7
$ws.Range("X5").Value = 1

# Row 6
$ws.Range("L6").Value = 4
$ws.Range("M6").Value = 6
$ws.Range("N6").Value = 7
$ws.Range("O6").Value = 8
$ws.Range("P6").Value = 9
$ws.Range("Q6").Value = 9
$ws.Range("R6").Value = 9
$ws.Range("S6").Value = 9
$ws.Range("T6").Value = 9
$ws.Range("U6").Value = 8
$ws.Range("V6").Value = 7
$ws.Range("W6").Value = 6
$ws.Range("X6").Value = 2

# Row 7
$ws.Range("L7").Value = 3
$ws.Range("M7").Value = 5
$ws.Range("N7").Value = 7
$ws.Range("O7").Value = 7
$ws.Range("Q7").Value = 9
$ws.Range("R7").Value = 9
$ws.Range("T7").Value = 8
$ws.Range("U7").Value = 8
$ws.Range("X7").Value = 1

# Row 8
$ws.Range("L8").Value = 3
$ws.Range("M8").Value = 3
$ws.Range("N8").Value = 5
$ws.Range("O8").Value = 6
$ws.Range("P8").Value = 7
$ws.Range("Q8").Value = 7
$ws.Range("S8").Value = 7
$ws.Range("T8").Value = 8
$ws.Range("U8").Value = 8
$ws.Range("V8").Value = 7
$ws.Range("W8").Value = 6
$ws.Range("X8").Value = 1

# Row 9
$ws.Range("N9").Value = 6
$ws.Range("O9").Value = 7
$ws.Range("P9").Value = 7
$ws.Range("Q9").Value = 7
$ws.Range("R9").Value = 7
$ws.Range("S9").Value = 7
$ws.Range("T9").Value = 7
$ws.Range("U9").Value = 7
$ws.Range("V9").Value = 7
$ws.Range("W9").Value = 6
$ws.Range("X9").Value = 2

# Row 10
$ws.Range("M10").Value = 4
$ws.Range("O10").Value = 7
$ws.Range("P10").Value = 7
$ws.Range("Q10").Value = 7
$ws.Range("R10").Value = 8
$ws.Range("S10").Value = 7
$ws.Range("T10").Value = 7
$ws.Range("U10").Value = 7
$ws.Range("V10").Value = 7
$ws.Range("W10").Value = 5
$ws.Range("X10").Value = 2

# Update selection to match the new active cell (X16)
$ws.Range("X16").Select()
